$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 8887
$ws.Range("E2").Value = 115
$ws.Range("F2").Value = 115
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = -6
$ws.Range("I2").Value = -4
$ws.Range("J2").Value = -2
$ws.Range("K2").Value = 5113
$ws.Range("L2").Value = 2477
$ws.Range("M2").Value = 2636
$ws.Range("N2").Value = 2627
$ws.Range("O2").Value = 9
$ws.Range("P2").Value = 208
$ws.Range("Q2").Value = 201
$ws.Range("R2").Value = 38
$ws.Range("S2").Value = -360
$ws.Range("T2").Value = 39
$ws.Range("U2").Value = 162
$ws.Range("V2").Value = 1195
$ws.Range("W2").Value = 1.29
$ws.Range("X2").Value = -0.06
$ws.Range("Y2").Value = -0.13
$ws.Range("Z2").Value = -0.11
$ws.Range("AA2").Value = 93.97
$ws.Range("AB2").Value = 1177.44
$ws.Range("AC2").Value = -84
$ws.Range("AD2").Value = -290.49
$ws.Range("AE2").Value = 66896
$ws.Range("AF2").Value = 0.37
$ws.Range("AG2").Value = 600
$ws.Range("AH2").Value = 2.45
$ws.Range("AI2").Value = -672.92
$ws.Range("AJ2").Value = 4160347
$ws.Range("D3").Value = 7810
$ws.Range("E3").Value = 99
$ws.Range("F3").Value = 99
$ws.Range("G3").Value = 32
$ws.Range("H3").Value = 21
$ws.Range("I3").Value = 20
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 5205
$ws.Range("L3").Value = 2579
$ws.Range("M3").Value = 2626
$ws.Range("N3").Value = 2616
$ws.Range("O3").Value = 10
$ws.Range("P3").Value = 208
$ws.Range("Q3").Value = 364
$ws.Range("R3").Value = -219
$ws.Range("S3").Value = -10
$ws.Range("T3").Value = 29
$ws.Range("U3").Value = 335
$ws.Range("V3").Value = 1263
$ws.Range("W3").Value = 1.26
$ws.Range("X3").Value = 0.27
$ws.Range("Y3").Value = 0.77
$ws.Range("Z3").Value = 0.41
$ws.Range("AA3").Value = 98.22
$ws.Range("AB3").Value = 1169.3
$ws.Range("AC3").Value = 485
$ws.Range("AD3").Value = 42.91
$ws.Range("AE3").Value = 66610
$ws.Range("AF3").Value = 0.31
$ws.Range("AG3").Value = 600
$ws.Range("AH3").Value = 2.88
$ws.Range("AI3").Value = 116.85
$ws.Range("AJ3").Value = 4160347
$ws.Range("D4").Value = 7494
$ws.Range("E4").Value = 154
$ws.Range("F4").Value = 168
$ws.Range("G4").Value = 24
$ws.Range("H4").Value = 29
$ws.Range("I4").Value = 27
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 5185
$ws.Range("L4").Value = 2545
$ws.Range("M4").Value = 2640
$ws.Range("N4").Value = 2628
$ws.Range("O4").Value = 11
$ws.Range("P4").Value = 208
$ws.Range("Q4").Value = 377
$ws.Range("R4").Value = 163
$ws.Range("S4").Value = -503
$ws.Range("T4").Value = 47
$ws.Range("U4").Value = 331
$ws.Range("V4").Value = 829
$ws.Range("W4").Value = 2.05
$ws.Range("X4").Value = 0.39
$ws.Range("Y4").Value = 1.04
$ws.Range("Z4").Value = 0.56
$ws.Range("AA4").Value = 96.39
$ws.Range("AB4").Value = 1173.03
$ws.Range("AC4").Value = 657
$ws.Range("AD4").Value = 33.8
$ws.Range("AE4").Value = 66929
$ws.Range("AF4").Value = 0.33
$ws.Range("AG4").Value = 600
$ws.Range("AH4").Value = 2.7
$ws.Range("AI4").Value = 86.24
$ws.Range("AJ4").Value = 4160347
$ws.Range("D5").Value = 8370
$ws.Range("E5").Value = 79
$ws.Range("F5").Value = 79
$ws.Range("G5").Value = 16
$ws.Range("H5").Value = 15
$ws.Range("I5").Value = 14
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 5390
$ws.Range("L5").Value = 2788
$ws.Range("M5").Value = 2601
$ws.Range("N5").Value = 2589
$ws.Range("O5").Value = 12
$ws.Range("P5").Value = 208
$ws.Range("Q5").Value = -152
$ws.Range("R5").Value = -210
$ws.Range("S5").Value = 183
$ws.Range("T5").Value = 60
$ws.Range("U5").Value = -212
$ws.Range("V5").Value = 1068
$ws.Range("W5").Value = 0.94
$ws.Range("X5").Value = 0.17
$ws.Range("Y5").Value = 0.54
$ws.Range("Z5").Value = 0.27
$ws.Range("AA5").Value = 107.2
$ws.Range("AB5").Value = 1165.84
$ws.Range("AC5").Value = 336
$ws.Range("AD5").Value = 68.28
$ws.Range("AE5").Value = 65929
$ws.Range("AF5").Value = 0.35
$ws.Range("AG5").Value = 600
$ws.Range("AH5").Value = 2.61
$ws.Range("AI5").Value = 168.5
$ws.Range("AJ5").Value = 4160347
$ws.Range("D6").Value = 8407
$ws.Range("E6").Value = 52
$ws.Range("F6").Value = 52
$ws.Range("G6").Value = 13
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("K6").Value = 5265
$ws.Range("L6").Value = 2661
$ws.Range("M6").Value = 2604
$ws.Range("N6").Value = 2592
$ws.Range("P6").Value = 208
$ws.Range("Q6").Value = 276
$ws.Range("R6").Value = -55
$ws.Range("S6").Value = -142
$ws.Range("T6").Value = 30
$ws.Range("U6").Value = 246
$ws.Range("V6").Value = 976
$ws.Range("W6").Value = 0.62
$ws.Range("X6").Value = 0.01
$ws.Range("Y6").Value = 0.04
$ws.Range("Z6").Value = 0.02
$ws.Range("AA6").Value = 102.18
$ws.Range("AB6").Value = 1155.47
$ws.Range("AC6").Value = 23
$ws.Range("AD6").Value = 747.91
$ws.Range("AE6").Value = 66000
$ws.Range("AF6").Value = 0.26
$ws.Range("AG6").Value = 600
$ws.Range("AH6").Value = 3.53
$ws.Range("AI6").Value = 2491.76
$ws.Range("AJ6").Value = 4160347

$ws.Range("D7:AJ9").ClearContents()
